$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up extra "alternative path/filename" columns that were scratch data ---
# Row 1: extra styled-but-empty cells D1:H1 (only A1:C1 remain)
$ws.Range("D1:H1").Clear()
# Row 2: extra alternative parameter-path cells
$ws.Range("C2:E2").Clear()
# Row 6: extra alternative parameter-filename cells
$ws.Range("C6:E6").Clear()
$ws.Range("I6").Clear()
$ws.Range("M6").Clear()
$ws.Range("S6").Clear()
$ws.Range("AB6").Clear()
# Row 8: extra output-data-filename cell
$ws.Range("C8").Clear()
# Row 9: extra ospm-path cell
$ws.Range("D9").Clear()

# --- Update the header of column C from "Alternatives" to "Comments" ---
$ws.Range("C1").Value = "Comments"

# Re-apply the "Neutral" cell style to C1 (same look as before, without the old
# unused "Bad" style) and drop the now-unused "Bad" cell style definition.
$ws.Range("C1").Style = "Neutral"
$wb.Styles.Item("Bad").Delete()

